# RPA datasets push 2024-04-09
# Insert two new IPO rows at the top of the table (below the header) and
# drop the two oldest rows from the bottom, keeping the table at 20 data
# rows (rows 2-21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row right below the header for the newest IPO (하스),
# pushing the former row 2 (노브랜드) down to row 3. Clear the formatting
# Excel copies down from the bold header row so the new row stays plain,
# like every other data row.
$ws.Rows("2:2").Insert()
$ws.Rows("2:2").ClearFormats()
$ws.Range("A2").Value = "하스"
$ws.Range("B2").Value = "2024.05.16~05.22"
$ws.Range("C2").Value = "9,000~12,000"
$ws.Range("D2").Value = "-"
$ws.Range("E2").Value = 16290
$ws.Range("F2").Value = "삼성증권"

# Insert a second new row after 노브랜드 for KB스팩28호, pushing the
# remaining former rows (starting at 아이씨티케이) down by one more.
$ws.Rows("4:4").Insert()
$ws.Range("A4").Value = "KB스팩28호"
$ws.Range("B4").Value = "2024.04.29~04.30"
$ws.Range("C4").Value = "2,000~2,000"
$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = 10000
$ws.Range("F4").Value = "KB증권"

# The table keeps a fixed size of 20 data rows, so drop the two oldest
# entries that were pushed past row 21 (하나스팩31호, 케이엔알시스템).
$ws.Rows("22:23").Delete()
